$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 9975.75
$ws.Range("I86").Value = 9966.333000000001
$ws.Range("K86").Value = 9966.333000000001
$ws.Range("M86").Value = -8843.333000000001
$ws.Range("H89").Value = 9975.75
$ws.Range("I89").Value = 9966.333000000001
$ws.Range("K89").Value = 49831.665
$ws.Range("M89").Value = -44215.665
$ws.Range("H132").Value = 1923.138
$ws.Range("I132").Value = 1404.0869
$ws.Range("K132").Value = 4212.2607
$ws.Range("M132").Value = -1682.2607
$ws.Range("H137").Value = 1607.6719
$ws.Range("I137").Value = 1520.4857
$ws.Range("K137").Value = 4561.4571
$ws.Range("M137").Value = -2011.4571
$ws.Range("H138").Value = 2501.0425
$ws.Range("I138").Value = 2381.158
$ws.Range("K138").Value = 7143.474
$ws.Range("M138").Value = -2003.474
$ws.Range("H141").Value = 3109.9167
$ws.Range("I141").Value = 3012.842
$ws.Range("J141").Value = 3478.8
$ws.Range("K141").Value = 9038.526
$ws.Range("L141").Value = 10436.4
$ws.Range("M141").Value = -3858.526
$ws.Range("N141").Value = -20796.4

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2919.1936
$ws.Range("I2").Value = 2675.52
$ws.Range("J2").Value = 3934.5
$ws.Range("K2").Value = 2675.52
$ws.Range("L2").Value = 3934.5
$ws.Range("M2").Value = -2562.52
$ws.Range("N2").Value = -4160.5
$ws.Range("H116").Value = 2919.1936
$ws.Range("I116").Value = 2675.52
$ws.Range("J116").Value = 3934.5
$ws.Range("K116").Value = 2675.52
$ws.Range("L116").Value = 3934.5
$ws.Range("M116").Value = -381.52
$ws.Range("N116").Value = -8522.5
$ws.Range("H132").Value = 3411.9312
$ws.Range("I132").Value = 3037.0454
$ws.Range("K132").Value = 9111.136200000001
$ws.Range("M132").Value = -6581.136200000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2919.1936
$ws.Range("I3").Value = 2675.52
$ws.Range("J3").Value = 3934.5
$ws.Range("K3").Value = 2675.52
$ws.Range("L3").Value = 3934.5
$ws.Range("M3").Value = -2561.52
$ws.Range("N3").Value = -4162.5
$ws.Range("H20").Value = 12529.613
$ws.Range("I20").Value = 14718.409
$ws.Range("K20").Value = 14718.409
$ws.Range("M20").Value = -14471.409
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H76").Value = 12625
$ws.Range("J76").Value = 12625
$ws.Range("L76").Value = 12625
$ws.Range("N76").Value = -13255
$ws.Range("H79").Value = 12625
$ws.Range("J79").Value = 12625
$ws.Range("L79").Value = 12625
$ws.Range("N79").Value = -14809
$ws.Range("H105").Value = 3032.4482
$ws.Range("I105").Value = 3120.0908
$ws.Range("J105").Value = 2978.889
$ws.Range("K105").Value = 3120.0908
$ws.Range("L105").Value = 2978.889
$ws.Range("M105").Value = -1373.0908
$ws.Range("N105").Value = -6472.889
$ws.Range("H107").Value = 2266.5
$ws.Range("I107").Value = 2154.3635
$ws.Range("J107").Value = 3500
$ws.Range("K107").Value = 2154.3635
$ws.Range("L107").Value = 3500
$ws.Range("M107").Value = -234.3634999999999
$ws.Range("N107").Value = -7340
$ws.Range("H134").Value = 5891.974
$ws.Range("I134").Value = 2711.1143
$ws.Range("J134").Value = 33724.5
$ws.Range("K134").Value = 8133.342900000001
$ws.Range("L134").Value = 101173.5
$ws.Range("M134").Value = -5598.342900000001
$ws.Range("N134").Value = -106243.5
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7598.8887
$ws.Range("I99").Value = 7778.1
$ws.Range("J99").Value = 7374.875
$ws.Range("K99").Value = 7778.1
$ws.Range("L99").Value = 7374.875
$ws.Range("M99").Value = -6280.1
$ws.Range("N99").Value = -10370.875
$ws.Range("H107").Value = 2296.3635
$ws.Range("J107").Value = 2615.75
$ws.Range("L107").Value = 2615.75
$ws.Range("N107").Value = -6455.75
$ws.Range("H126").Value = 7598.8887
$ws.Range("I126").Value = 7778.1
$ws.Range("J126").Value = 7374.875
$ws.Range("K126").Value = 23334.3
$ws.Range("L126").Value = 22124.625
$ws.Range("M126").Value = -20864.3
$ws.Range("N126").Value = -27064.625
$ws.Range("H132").Value = 6709.727
$ws.Range("I132").Value = 7113.5713
$ws.Range("J132").Value = 6003
$ws.Range("K132").Value = 21340.7139
$ws.Range("L132").Value = 18009
$ws.Range("M132").Value = -18810.7139
$ws.Range("N132").Value = -23069

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1120.2094
$ws.Range("I5").Value = 812.95654
$ws.Range("J5").Value = 1473.55
$ws.Range("K5").Value = 2438.86962
$ws.Range("L5").Value = 4420.65
$ws.Range("M5").Value = -2326.86962
$ws.Range("N5").Value = -4644.65
$ws.Range("H80").Value = 50000
$ws.Range("I80").Value = 50000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 150000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -149064
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 50000
$ws.Range("I83").Value = 50000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 450000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -445320
$ws.Range("N83").ClearContents()
$ws.Range("H114").Value = 2093.9092
$ws.Range("I114").Value = 249.33333
$ws.Range("J114").Value = 2785.625
$ws.Range("K114").Value = 747.99999
$ws.Range("L114").Value = 8356.875
$ws.Range("M114").Value = 2506.00001
$ws.Range("N114").Value = -14864.875
$ws.Range("H115").Value = 1249.5
$ws.Range("I115").Value = 1249.5
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 3748.5
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -2573.5
$ws.Range("N115").ClearContents()
$ws.Range("H135").Value = 1120.2094
$ws.Range("I135").Value = 812.95654
$ws.Range("J135").Value = 1473.55
$ws.Range("K135").Value = 7316.60886
$ws.Range("L135").Value = 13261.95
$ws.Range("M135").Value = -4781.60886
$ws.Range("N135").Value = -18331.95

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8473.695
$ws.Range("I70").Value = 7849.364
$ws.Range("K70").Value = 7849.364
$ws.Range("M70").Value = -7579.364
$ws.Range("H73").Value = 8473.695
$ws.Range("I73").Value = 7849.364
$ws.Range("K73").Value = 7849.364
$ws.Range("M73").Value = -6913.364
$ws.Range("H80").Value = 9373
$ws.Range("I80").Value = 5071.6
$ws.Range("J80").Value = 14749.75
$ws.Range("K80").Value = 5071.6
$ws.Range("L80").Value = 14749.75
$ws.Range("M80").Value = -4073.6
$ws.Range("N80").Value = -16745.75
$ws.Range("H83").Value = 9373
$ws.Range("I83").Value = 5071.6
$ws.Range("J83").Value = 14749.75
$ws.Range("K83").Value = 25358
$ws.Range("L83").Value = 73748.75
$ws.Range("M83").Value = -20366
$ws.Range("N83").Value = -83732.75
$ws.Range("H132").Value = 14362.5
$ws.Range("I132").Value = 16396
$ws.Range("J132").Value = 6228.5
$ws.Range("K132").Value = 49188
$ws.Range("L132").Value = 18685.5
$ws.Range("M132").Value = -46658
$ws.Range("N132").Value = -23745.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 185.27272
$ws.Range("I55").Value = 139.5
$ws.Range("J55").Value = 240.2
$ws.Range("K55").Value = 139.5
$ws.Range("L55").Value = 240.2
$ws.Range("M55").Value = 33.5
$ws.Range("N55").Value = -586.2
$ws.Range("H69").Value = 1299749.6
$ws.Range("J69").Value = 1299749.6
$ws.Range("L69").Value = 1299749.6
$ws.Range("N69").Value = -1301371.6
$ws.Range("H72").Value = 1299749.6
$ws.Range("J72").Value = 1299749.6
$ws.Range("L72").Value = 3899248.8
$ws.Range("N72").Value = -3907360.8
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H132").Value = 3608.0557
$ws.Range("I132").Value = 2620.25
$ws.Range("K132").Value = 7860.75
$ws.Range("M132").Value = -5330.75
$ws.Range("H135").Value = 46899.6
$ws.Range("J135").Value = 46899.6
$ws.Range("L135").Value = 46899.6
$ws.Range("N135").Value = -57039.6
$ws.Range("H141").Value = 59997
$ws.Range("J141").Value = 59997
$ws.Range("L141").Value = 59997
$ws.Range("N141").Value = -70357
